$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A19").Value = "seriesId"
$ws.Range("B19").Value = "Client"
$ws.Range("C19").Value = "Client"
$ws.Range("D19").Value = "No"

$ws.Range("B20").Select()
